$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")
$ws.Columns.Item(8).Insert()
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(2, 8).Value = "stock"
$ws.Cells.Item(3, 8).Value = "stock"
$ws.Cells.Item(4, 8).Value = "stock"
